$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = 3.1
$ws.Range("R2").Value = 1.44
$ws.Range("AC2").Value = 6.5
$ws.Range("AO2").Value = 13
$ws.Range("AW2").Value = 5.5
$ws.Range("R3").Value = 1.44
$ws.Range("R4").Value = 1.7
$ws.Range("R5").Value = 1.48
$ws.Range("M7").Value = 1.07
$ws.Range("N7").Value = 9
$ws.Range("Q7").Value = 2.1
$ws.Range("R7").Value = 1.7
$ws.Range("J8").Value = 3
$ws.Range("Q8").Value = 2.15
$ws.Range("R8").Value = 1.67
$ws.Range("S8").Value = 1.5
$ws.Range("T8").Value = 2.5
$ws.Range("W8").Value = 6.5
$ws.Range("Z8").Value = 21
$ws.Range("AG8").Value = 8.5
$ws.Range("AO8").Value = 13
$ws.Range("AP8").Value = 26
$ws.Range("AT8").Value = 2.5
$ws.Range("AV8").Value = 67
$ws.Range("G9").Value = 1.6
$ws.Range("H9").Value = 3.75
$ws.Range("I9").Value = 6.25
$ws.Range("J9").Value = 2.2
$ws.Range("U9").Value = 1.95
$ws.Range("V9").Value = 1.8
$ws.Range("AE9").Value = 17
$ws.Range("AH9").Value = 29
$ws.Range("AI9").Value = 19
$ws.Range("AK9").Value = 41
$ws.Range("AN9").Value = 3.5
$ws.Range("AW9").Value = 7
$ws.Range("AX9").Value = 29
$ws.Range("BB9").Value = 301
$ws.Range("G10").Value = 2.6
$ws.Range("K10").Value = 1.95
$ws.Range("S10").Value = 1.53
$ws.Range("T10").Value = 2.38
$ws.Range("U10").Value = 2
$ws.Range("V10").Value = 1.75
$ws.Range("W10").Value = 7
$ws.Range("AA10").Value = 23
$ws.Range("AC10").Value = 7
$ws.Range("AD10").Value = 5.5
$ws.Range("AJ10").Value = 29
$ws.Range("AM10").Value = 451
$ws.Range("AS10").Value = 251
$ws.Range("AT10").Value = 2.38
$ws.Range("AY10").Value = 29
$ws.Range("M12").Value = 1.04
$ws.Range("N12").Value = 13
$ws.Range("G14").Value = 1.75
$ws.Range("I14").Value = 5
$ws.Range("J14").Value = 2.38
$ws.Range("AK14").Value = 41
$ws.Range("G20").Value = 30
$ws.Range("H20").Value = 7.8
$ws.Range("J20").Value = 20
$ws.Range("L20").Value = 1.32
$ws.Range("R20").Value = 3.4
$ws.Range("X20").Value = 700
$ws.Range("Y20").Value = 150
$ws.Range("AA20").Value = 800
$ws.Range("AB20").Value = 450
$ws.Range("AC20").Value = 23
$ws.Range("AF20").Value = 250
$ws.Range("AG20").Value = 11.5
$ws.Range("AN20").Value = 27
$ws.Range("AO20").Value = 250
$ws.Range("AP20").Value = 120
$ws.Range("AU20").Value = 12
$ws.Range("AV20").Value = 110
$ws.Range("AW20").Value = 3.2
$ws.Range("AZ20").Value = 7.3
$ws.Range("BA20").Value = 27
